$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the changed Price/Volume cells to remain text so Excel does not
# auto-convert numeric- or percent-looking strings into Number/Percentage cells.
$textCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D20", "D21", "D23", "D24", "D25", "D26", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D47", "D48", "D50", "D51", "E2", "E3", "E4", "E5", "E6", "E7", "E8", "E9", "E10", "E11", "E12", "E13", "E14", "E15", "E16", "E17", "E18", "E20", "E21", "E22", "E23", "E24", "E25", "E26", "E27", "E39", "E40", "E41", "E42", "E43", "E44", "E45", "E46", "E47", "E48", "E50", "E51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "328.07"
$ws.Range("E2").Value = "-0.76%"
$ws.Range("D3").Value = "44.10"
$ws.Range("E3").Value = "6.09%"
$ws.Range("D4").Value = "5.477"
$ws.Range("E4").Value = "-3.88%"
$ws.Range("D5").Value = "0.08084"
$ws.Range("E5").Value = "-3.70%"
$ws.Range("E6").Value = "-1.40%"
$ws.Range("D7").Value = "4.300"
$ws.Range("E7").Value = "-3.96%"
$ws.Range("D8").Value = "1.891"
$ws.Range("E8").Value = "-5.47%"
$ws.Range("D9").Value = "2.745"
$ws.Range("E9").Value = "-6.58%"
$ws.Range("D10").Value = "0.9396"
$ws.Range("E10").Value = "1.81%"
$ws.Range("D11").Value = "0.1223"
$ws.Range("E11").Value = "-4.47%"
$ws.Range("D12").Value = "0.1894"
$ws.Range("E12").Value = "-4.02%"
$ws.Range("D13").Value = "0.09699"
$ws.Range("E13").Value = "2.30%"
$ws.Range("D14").Value = "0.04131"
$ws.Range("E14").Value = "6.36%"
$ws.Range("D15").Value = "0.1069"
$ws.Range("E15").Value = "0.72%"
$ws.Range("D16").Value = "0.001276"
$ws.Range("E16").Value = "-2.02%"
$ws.Range("D17").Value = "0.006035"
$ws.Range("E17").Value = "-1.21%"
$ws.Range("D18").Value = "3.570"
$ws.Range("E18").Value = "4.25%"
$ws.Range("D20").Value = "8.513"
$ws.Range("E20").Value = "-4.53%"
$ws.Range("D21").Value = "0.1349"
$ws.Range("E21").Value = "-0.99%"
$ws.Range("E22").Value = "-0.62%"
$ws.Range("D23").Value = "0.04371"
$ws.Range("E23").Value = "-0.71%"
$ws.Range("D24").Value = "0.001235"
$ws.Range("E24").Value = "-3.06%"
$ws.Range("D25").Value = "0.004301"
$ws.Range("E25").Value = "-1.26%"
$ws.Range("D26").Value = "0.0001234"
$ws.Range("E26").Value = "3.61%"
$ws.Range("E27").Value = "0.25%"
$ws.Range("D39").Value = "0.02649"
$ws.Range("E39").Value = "-6.71%"
$ws.Range("D40").Value = "0.05458"
$ws.Range("E40").Value = "-1.08%"
$ws.Range("D41").Value = "0.007685"
$ws.Range("E41").Value = "-3.36%"
$ws.Range("D42").Value = "0.009719"
$ws.Range("E42").Value = "8.24%"
$ws.Range("D43").Value = "0.1389"
$ws.Range("E43").Value = "-3.08%"
$ws.Range("D44").Value = "0.002126"
$ws.Range("E44").Value = "1.62%"
$ws.Range("D45").Value = "0.009896"
$ws.Range("E45").Value = "-15.27%"
$ws.Range("E46").Value = "2.56%"
$ws.Range("D47").Value = "0.00000000752"
$ws.Range("E47").Value = "0.25%"
$ws.Range("D48").Value = "0.003553"
$ws.Range("E48").Value = "2.60%"
$ws.Range("D50").Value = "0.00002107"
$ws.Range("E50").Value = "0.25%"
$ws.Range("D51").Value = "0.0002006"
$ws.Range("E51").Value = "0.25%"
